# Blogspot links & pictures added
# Adds two new backlink rows (16 and 17) to Sheet1, each with a blogspot
# URL in column B, an email address in column C (as a mailto hyperlink)
# and the literal "realbreeze" in column D - matching the existing
# pattern used by rows 10-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# NOTE: the shared-string table records strings in the order their text is
# first written, so row 17's URL is entered before row 16's URL to land on
# the expected shared-string indices (16 = surigaoparadise, 17 =
# boholtraveladventures).
$ws.Range("B17").Value = "https://surigaoparadise.blogspot.com/"
$ws.Range("B16").Value = "https://boholtraveladventures.blogspot.com/"

$ws.Range("C16").Value = "realbreezemark@gmail.com"
$ws.Range("D16").Value = "realbreeze"

$ws.Range("C17").Value = "realbreezemark@gmail.com"
$ws.Range("D17").Value = "realbreeze"

# Wire up the hyperlinks (mailto: for column C, the blog URL itself for B17).
$ws.Hyperlinks.Add($ws.Range("C16"), "mailto:realbreezemark@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B17"), "https://surigaoparadise.blogspot.com/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C17"), "mailto:realbreezemark@gmail.com") | Out-Null

# Adding a hyperlink re-applies direct formatting; restore the shared
# "Hyperlink" cell style so these cells match the existing linked cells
# (B2:B8, C10:C15) instead of keeping a separate direct-format copy.
$ws.Range("C16").Style = "Hyperlink"
$ws.Range("B17").Style = "Hyperlink"
$ws.Range("C17").Style = "Hyperlink"
